$d = $word.ActiveDocument

$old = "Dataset, consisting of 147 unique instances, was distributed as part of the master course on Data Mining in Bioinformatics at the Faculty of Mathematics, University of Belgrade"
$new = "Dataset, consisting of 147 unique instances, was obtained from the Pacemaker Center of the Clinical Center of Serbia and pertains to patients with electrical cardioversion performed from 2014 to 2019"

$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

$style = $d.Styles("Default Paragraph Font")
$style.Hidden = $false
